$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("despesas")
$ws.Name = "Folha1"
$ws.Activate()
$ws.Range("H20").Select()
